# 2.c.1.1b.xlsx - add the "2023" column (T) to the CPI table.
#
# The sheet has a year header row (row 4, columns D..S = 2007..2022) and
# ten data rows (rows 5..14) below it. We extend the table one column to
# the right (column T) with the 2023 figures, matching the look & feel of
# the existing 2022 column (S) - same fonts/borders, but with a "0.0"
# number format instead of "General"/"#,##0.0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header cell T4 : 2023, same style as S4 -------------------------------
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value2 = 2023

# --- T5 ("Kyrgyz Republic" bold total row), style like S5 but "0.0" -------
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value2 = 103.3752971770694
$ws.Range("T5").NumberFormat = "0.0"

# --- T6:T13 (regular oblast rows), style like S6 but "0.0" -----------------
$values = @{
  6  = 104.28289267885566
  7  = 105.54952929965596
  8  = 103.51673393645765
  9  = 106.47711005273266
  10 = 107.06485907396235
  11 = 109.45786616400459
  12 = 101.11876601355125
  13 = 102.74015470892634
}
foreach ($r in 6..13) {
  $ws.Range("S6").Copy()
  $ws.Range("T$r").PasteSpecial(-4122)
  $ws.Range("T$r").Value2 = $values[$r]
  $ws.Range("T$r").NumberFormat = "0.0"
}

# --- T14 (bottom bordered row), style like S14 but "0.0" -------------------
$ws.Range("S14").Copy()
$ws.Range("T14").PasteSpecial(-4122)
$ws.Range("T14").Value2 = 103.21772010523679
$ws.Range("T14").NumberFormat = "0.0"

# --- row heights for the data block grow slightly to fit the new column ---
foreach ($r in 5..14) {
  $ws.Rows.Item($r).RowHeight = 14.25
}

# --- columns D:T get a touch narrower to make room for the extra column ---
$ws.Range("D1:T1").ColumnWidth = 7.711495535714286

# --- return cursor to the top-left corner ----------------------------------
$ws.Range("A1").Select()
